$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', StandardScaler()),`n                ('selector',`n                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),`n                ('model',`n                 LogisticRegression(C=1, max_iter=1000, random_state=42,`n                                    solver='liblinear'))])"
$ws.Range("B2").Value = 0.6952380952380952
$ws.Range("C2").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': StandardScaler(), 'model__solver': 'liblinear', 'model__penalty': 'l2', 'model__class_weight': None, 'model__C': 1}"
$ws.Range("D2").Value = 0.3333333333333333
$ws.Range("E2").Value = "[1 1 0 0 1 0 0 0 0 1 0 1]"
$ws.Range("F2").Value = "[0 0 1 0 0 1 1 0 1 1 1 1]"
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 0.6736607142857142
$ws.Range("I2").Value = 0.02294307765069291
$ws.Range("J2").Value = 0.5801587301587302
$ws.Range("K2").Value = 0.05163003293260471

# Row 3
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',`n                                                     random_state=42))),`n                ('model',`n                 LogisticRegression(C=0.001, class_weight='balanced',`n                                    max_iter=1000, penalty='l1',`n                                    random_state=42, solver='saga'))])"
$ws.Range("B3").Value = 0.6285714285714284
$ws.Range("C3").Value = "{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': None, 'model__solver': 'saga', 'model__penalty': 'l1', 'model__class_weight': 'balanced', 'model__C': 0.001}"
$ws.Range("D3").Value = 0.7692307692307692
$ws.Range("E3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0]"
$ws.Range("F3").Value = "[0 1 0 1 1 0 1 0 0 1 1 0]"
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.6616071428571429
$ws.Range("I3").Value = 0.02780134286098932
$ws.Range("J3").Value = 0.5494897959183673
$ws.Range("K3").Value = 0.0768380854278453

# Row 4
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),`n                ('model',`n                 LogisticRegression(C=0.0001, class_weight='balanced',`n                                    max_iter=1000, random_state=42,`n                                    solver='saga'))])"
$ws.Range("B4").Value = 0.6380952380952382
$ws.Range("C4").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__solver': 'saga', 'model__penalty': 'l2', 'model__class_weight': 'balanced', 'model__C': 0.0001}"
$ws.Range("D4").Value = 0.5714285714285715
$ws.Range("E4").Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Range("F4").Value = "[0 0 0 1 1 1 1 0 1 0 0 1]"
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.6779503105590062
$ws.Range("I4").Value = 0.03239758790717629
$ws.Range("J4").Value = 0.5616977225672877
$ws.Range("K4").Value = 0.08836305531615225

# Row 5
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',`n                                                     random_state=42))),`n                ('model',`n                 LogisticRegression(C=0.001, class_weight='balanced',`n                                    max_iter=1000, penalty='l1',`n                                    random_state=42, solver='saga'))])"
$ws.Range("B5").Value = 0.6190476190476191
$ws.Range("C5").Value = "{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': None, 'model__solver': 'saga', 'model__penalty': 'l1', 'model__class_weight': 'balanced', 'model__C': 0.001}"
$ws.Range("D5").Value = 0.7999999999999999
$ws.Range("E5").Value = "[1 1 0 0 0 0 1 0 1 1 1 1]"
$ws.Range("F5").Value = "[0 1 0 0 0 1 1 1 1 1 1 1]"
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 0.667340521114106
$ws.Range("I5").Value = 0.03268945269624195
$ws.Range("J5").Value = 0.5324348607367476
$ws.Range("K5").Value = 0.07788552347151424

# Row 6
$ws.Range("A6").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f3a6c5dd400>),`n                ('model',`n                 LogisticRegression(C=5, max_iter=1000, random_state=42,`n                                    solver='liblinear'))])"
$ws.Range("B6").Value = 0.619047619047619
$ws.Range("C6").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f3b000c7550>, 'scaler': None, 'model__solver': 'liblinear', 'model__penalty': 'l2', 'model__class_weight': None, 'model__C': 5}"
$ws.Range("D6").Value = 0.4615384615384615
$ws.Range("E6").Value = "[1 1 1 1 0 0 0 0 1 1 0 0]"
$ws.Range("F6").Value = "[1 1 0 0 0 1 1 1 1 0 1 0]"
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 0.6886684303350968
$ws.Range("I6").Value = 0.02529738266553202
$ws.Range("J6").Value = 0.5599647266313932
$ws.Range("K6").Value = 0.05978725111305908
